$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        $Worksheet,
        [int]$Row,
        [int]$Col,
        [string]$Text
    )
    $cell = $Worksheet.Cells.Item($Row, $Col)
    # Force the cell to text format so the numeric-looking string is preserved
    # verbatim (keeps leading/trailing zeros, decimal formatting, etc.)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    # Reset the style back to the workbook default so no stray formatting
    # is introduced on the cell.
    $cell.Style = "Normal"
}

# Row 88 (Grand Est / M)
Set-TextValue $ws 88 3 "121"
Set-TextValue $ws 88 5 "1023835.59"

# Row 130 (Hauts-de-France / I)
Set-TextValue $ws 130 3 "1153"
Set-TextValue $ws 130 5 "9575831.11"

# Row 157 (La Reunion / C)
Set-TextValue $ws 157 3 "18"
Set-TextValue $ws 157 4 "18"
Set-TextValue $ws 157 5 "45000.00"

# Row 170 (Martinique / C)
Set-TextValue $ws 170 3 "55"
Set-TextValue $ws 170 4 "53"
Set-TextValue $ws 170 5 "148174.00"

# Row 171 (Martinique / F)
Set-TextValue $ws 171 3 "54"
Set-TextValue $ws 171 4 "54"
Set-TextValue $ws 171 5 "144500.00"

# Row 172 (Martinique / G)
Set-TextValue $ws 172 3 "150"
Set-TextValue $ws 172 4 "148"
Set-TextValue $ws 172 5 "430493.00"

# Row 173 (Martinique / H)
Set-TextValue $ws 173 3 "37"
Set-TextValue $ws 173 4 "37"
Set-TextValue $ws 173 5 "90204.22"

# Row 174 (Martinique / I)
Set-TextValue $ws 174 3 "106"
Set-TextValue $ws 174 4 "105"
Set-TextValue $ws 174 5 "448831.34"

# Row 175 (Martinique / J)
Set-TextValue $ws 175 3 "8"
Set-TextValue $ws 175 4 "8"
Set-TextValue $ws 175 5 "17500.00"

# Row 178 (Martinique / M)
Set-TextValue $ws 178 3 "57"
Set-TextValue $ws 178 4 "55"
Set-TextValue $ws 178 5 "189394.00"

# Row 179 (Martinique / N)
Set-TextValue $ws 179 3 "50"
Set-TextValue $ws 179 4 "48"
Set-TextValue $ws 179 5 "128760.60"

# Row 180 (Martinique / P)
Set-TextValue $ws 180 3 "19"
Set-TextValue $ws 180 4 "19"
Set-TextValue $ws 180 5 "39500.00"

# Row 183 (Martinique / S)
Set-TextValue $ws 183 3 "78"
Set-TextValue $ws 183 4 "75"
Set-TextValue $ws 183 5 "223618.00"

Write-Host "Done applying 2020-12-14 Fonds de solidarite updates"
